$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13290.909
$ws.Range("I18").Value = 12620
$ws.Range("J18").Value = 20000
$ws.Range("K18").Value = 12620
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = -12336
$ws.Range("N18").Value = -20568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 1607.7693
$ws.Range("I36").Value = 1408.4166
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 1408.4166
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = -1062.4166
$ws.Range("N36").Value = -4692

$ws.Range("H61").Value = 1699.5
$ws.Range("I61").Value = 1699.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1699.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1487.5

$ws.Range("H136").Value = 1699.5
$ws.Range("I136").Value = 1699.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5098.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2548.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29629.285
$ws.Range("I86").Value = 1866
$ws.Range("J86").Value = 50451.75
$ws.Range("K86").Value = 1866
$ws.Range("L86").Value = 50451.75
$ws.Range("M86").Value = -743
$ws.Range("N86").Value = -52697.75

$ws.Range("H89").Value = 29629.285
$ws.Range("I89").Value = 1866
$ws.Range("J89").Value = 50451.75
$ws.Range("K89").Value = 9330
$ws.Range("L89").Value = 252258.75
$ws.Range("M89").Value = -3714
$ws.Range("N89").Value = -263490.75

$ws.Range("H94").Value = 2423.7144
$ws.Range("I94").Value = 1604.8889
$ws.Range("J94").Value = 7336.6665
$ws.Range("K94").Value = 1604.8889
$ws.Range("L94").Value = 7336.6665
$ws.Range("M94").Value = -1153.8889
$ws.Range("N94").Value = -8238.666499999999

$ws.Range("H99").Value = 3402.6
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 3004.3333
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 3004.3333
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -6000.3333

$ws.Range("H134").Value = 2968.625
$ws.Range("I134").Value = 3178.4285
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 9535.2855
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -7000.2855
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 5154.2
$ws.Range("I105").Value = 1648.8572
$ws.Range("J105").Value = 13333.333
$ws.Range("K105").Value = 1648.8572
$ws.Range("L105").Value = 13333.333
$ws.Range("M105").Value = 98.14280000000008
$ws.Range("N105").Value = -16827.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 614
$ws.Range("I12").Value = 2999
$ws.Range("J12").Value = 17.75
$ws.Range("K12").Value = 8997
$ws.Range("L12").Value = 53.25
$ws.Range("M12").Value = -8824
$ws.Range("N12").Value = -399.25

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H76").Value = 1000
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2617

$ws.Range("H79").Value = 1000
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1674

$ws.Range("H114").Value = 589
$ws.Range("I114").Value = 178
$ws.Range("J114").Value = 1000
$ws.Range("K114").Value = 534
$ws.Range("L114").Value = 3000
$ws.Range("M114").Value = 2720
$ws.Range("N114").Value = -9508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H80").Value = 5600
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 6500
$ws.Range("K80").Value = 3800
$ws.Range("L80").Value = 6500
$ws.Range("M80").Value = -2802
$ws.Range("N80").Value = -8496

$ws.Range("H83").Value = 5600
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 6500
$ws.Range("K83").Value = 19000
$ws.Range("L83").Value = 32500
$ws.Range("M83").Value = -14008
$ws.Range("N83").Value = -42484

$ws.Range("H122").Value = 1360
$ws.Range("I122").Value = 1360
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4080
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1630
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 9499
$ws.Range("I132").Value = 8500
$ws.Range("J132").Value = 9998.5
$ws.Range("K132").Value = 25500
$ws.Range("L132").Value = 29995.5
$ws.Range("M132").Value = -22970
$ws.Range("N132").Value = -35055.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4980.846
$ws.Range("I7").Value = 4979.25
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 4979.25
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4867.25
$ws.Range("N7").Value = -5224

$ws.Range("H46").Value = 4027.6
$ws.Range("I46").Value = 1624.75
$ws.Range("J46").Value = 4901.364
$ws.Range("K46").Value = 1624.75
$ws.Range("L46").Value = 4901.364
$ws.Range("M46").Value = -1436.75
$ws.Range("N46").Value = -5277.364

$ws.Range("H55").Value = 2413.6667
$ws.Range("I55").Value = 966.3333
$ws.Range("J55").Value = 4584.6665
$ws.Range("K55").Value = 966.3333
$ws.Range("L55").Value = 4584.6665
$ws.Range("M55").Value = -793.3333
$ws.Range("N55").Value = -4930.6665

$ws.Range("H82").Value = 1952.6666
$ws.Range("I82").Value = 1617.7858
$ws.Range("J82").Value = 3124.75
$ws.Range("K82").Value = 1617.7858
$ws.Range("L82").Value = 3124.75
$ws.Range("M82").Value = -1256.7858
$ws.Range("N82").Value = -3846.75

$ws.Range("H85").Value = 1952.6666
$ws.Range("I85").Value = 1617.7858
$ws.Range("J85").Value = 3124.75
$ws.Range("K85").Value = 1617.7858
$ws.Range("L85").Value = 3124.75
$ws.Range("M85").Value = -369.7858000000001
$ws.Range("N85").Value = -5620.75

$ws.Range("H126").Value = 4980.846
$ws.Range("I126").Value = 4979.25
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 14937.75
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -12467.75
$ws.Range("N126").Value = -19940

$ws.Range("H136").Value = 5465.8125
$ws.Range("I136").Value = 4727.154
$ws.Range("J136").Value = 8666.666999999999
$ws.Range("K136").Value = 14181.462
$ws.Range("L136").Value = 26000.001
$ws.Range("M136").Value = -11631.462
$ws.Range("N136").Value = -31100.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H81").Value = 1998.5
$ws.Range("I81").Value = 1998.3334
$ws.Range("J81").Value = 1999
$ws.Range("K81").Value = 3996.6668
$ws.Range("L81").Value = 3998
$ws.Range("M81").Value = -2935.6668
$ws.Range("N81").Value = -6120

$ws.Range("H84").Value = 1998.5
$ws.Range("I84").Value = 1998.3334
$ws.Range("J84").Value = 1999
$ws.Range("K84").Value = 19983.334
$ws.Range("L84").Value = 19990
$ws.Range("M84").Value = -14679.334
$ws.Range("N84").Value = -30598

$ws.Range("H122").Value = 3154.8333
$ws.Range("I122").Value = 3425.8
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 10277.4
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -7827.400000000001
$ws.Range("N122").Value = -10300

$ws.Range("H132").Value = 2462.5833
$ws.Range("I132").Value = 2472.375
$ws.Range("J132").Value = 2443
$ws.Range("K132").Value = 7417.125
$ws.Range("L132").Value = 7329
$ws.Range("M132").Value = -4887.125
$ws.Range("N132").Value = -12389
